$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2042253521126761
$ws.Range("C2").Value = 0.528169014084507
$ws.Range("J2").Value = 0.01056338028169014
$ws.Range("P2").Value = 0.1514084507042253
$ws.Range("S2").Value = 0.1056338028169014
$ws.Range("B3").Value = 0.01851851851851852
$ws.Range("C3").Value = 0.04320987654320987
$ws.Range("J3").Value = 0.01234567901234568
$ws.Range("P3").Value = 0.6604938271604939
$ws.Range("S3").Value = 0.2654320987654321
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.05857740585774059
$ws.Range("D6").Value = 0.004184100418410041
$ws.Range("F6").Value = 0.05439330543933055
$ws.Range("J6").Value = 0.301255230125523
$ws.Range("O6").Value = 0.01673640167364017
$ws.Range("Q6").Value = 0.1213389121338912
$ws.Range("R6").Value = 0.06276150627615062
$ws.Range("S6").Value = 0.3807531380753138
$ws.Range("B7").Value = 0.1197916666666667
$ws.Range("D7").Value = 0.01041666666666667
$ws.Range("E7").Value = 0.005208333333333333
$ws.Range("F7").Value = 0.05208333333333334
$ws.Range("J7").Value = 0.1145833333333333
$ws.Range("O7").Value = 0.02083333333333333
$ws.Range("Q7").Value = 0.1770833333333333
$ws.Range("R7").Value = 0.06770833333333333
$ws.Range("S7").Value = 0.4322916666666667
$ws.Range("B8").Value = 0.08478260869565217
$ws.Range("D8").Value = 0.01521739130434783
$ws.Range("E8").Value = 0.002173913043478261
$ws.Range("F8").Value = 0.06304347826086956
$ws.Range("J8").Value = 0.1260869565217391
$ws.Range("Q8").Value = 0.1826086956521739
$ws.Range("R8").Value = 0.07173913043478261
$ws.Range("S8").Value = 0.4543478260869565
$ws.Range("B9").Value = 0.09836065573770492
$ws.Range("D9").Value = 0.01092896174863388
$ws.Range("F9").Value = 0.06557377049180328
$ws.Range("J9").Value = 0.1038251366120219
$ws.Range("O9").Value = 0.00546448087431694
$ws.Range("Q9").Value = 0.2021857923497268
$ws.Range("R9").Value = 0.08743169398907104
$ws.Range("S9").Value = 0.4262295081967213
$ws.Range("B10").Value = 0.1033057851239669
$ws.Range("D10").Value = 0.02809917355371901
$ws.Range("E10").Value = 0.002479338842975207
$ws.Range("F10").Value = 0.07851239669421488
$ws.Range("J10").Value = 0.1173553719008264
$ws.Range("O10").Value = 0.0115702479338843
$ws.Range("Q10").Value = 0.1867768595041322
$ws.Range("R10").Value = 0.08264462809917356
$ws.Range("S10").Value = 0.3892561983471075
$ws.Range("G11").Value = 0.1412103746397695
$ws.Range("J11").Value = 0.1095100864553314
$ws.Range("K11").Value = 0.2190201729106628
$ws.Range("L11").Value = 0.4956772334293948
$ws.Range("S11").Value = 0.0345821325648415
$ws.Range("G12").Value = 0.615819209039548
$ws.Range("J12").Value = 0.2824858757062147
$ws.Range("K12").Value = 0.01129943502824859
$ws.Range("L12").Value = 0.04519774011299435
$ws.Range("S12").Value = 0.04519774011299435
$ws.Range("J14").Value = 0.5
$ws.Range("S14").Value = 0.5
$ws.Range("F15").Value = 0.04324324324324325
$ws.Range("H15").Value = 0.1783783783783784
$ws.Range("I15").Value = 0.06486486486486487
$ws.Range("J15").Value = 0.3189189189189189
$ws.Range("K15").Value = 0.06486486486486487
$ws.Range("M15").Value = 0.01621621621621622
$ws.Range("O15").Value = 0.03243243243243243
$ws.Range("S15").Value = 0.2810810810810811
$ws.Range("F16").Value = 0.01630434782608696
$ws.Range("H16").Value = 0.1521739130434783
$ws.Range("I16").Value = 0.05978260869565218
$ws.Range("J16").Value = 0.3478260869565217
$ws.Range("K16").Value = 0.1304347826086956
$ws.Range("M16").Value = 0.02717391304347826
$ws.Range("N16").Value = 0.005434782608695652
$ws.Range("O16").Value = 0.04347826086956522
$ws.Range("S16").Value = 0.2173913043478261
$ws.Range("F17").Value = 0.0170316301703163
$ws.Range("H17").Value = 0.1873479318734793
$ws.Range("I17").Value = 0.0827250608272506
$ws.Range("J17").Value = 0.3625304136253041
$ws.Range("K17").Value = 0.0827250608272506
$ws.Range("M17").Value = 0.024330900243309
$ws.Range("N17").Value = 0.004866180048661801
$ws.Range("O17").Value = 0.04866180048661801
$ws.Range("S17").Value = 0.1897810218978102
$ws.Range("F18").Value = 0.02808988764044944
$ws.Range("H18").Value = 0.2584269662921349
$ws.Range("I18").Value = 0.08426966292134831
$ws.Range("J18").Value = 0.3370786516853932
$ws.Range("K18").Value = 0.1235955056179775
$ws.Range("M18").Value = 0.005617977528089887
$ws.Range("O18").Value = 0.02808988764044944
$ws.Range("S18").Value = 0.1348314606741573
$ws.Range("F19").Value = 0.02367688022284123
$ws.Range("H19").Value = 0.1942896935933148
$ws.Range("I19").Value = 0.0786908077994429
$ws.Range("J19").Value = 0.334958217270195
$ws.Range("K19").Value = 0.1225626740947075
$ws.Range("M19").Value = 0.02437325905292479
$ws.Range("O19").Value = 0.06337047353760446
$ws.Range("S19").Value = 0.1580779944289694
